# Update the "Test Cases" sheet (F2 -> Result = FAIL) and the
# "VerifyStartEndDateValidation" sheet (K2 -> Result = FAIL), reflecting
# that automated test runs now report a FAIL result instead of being blank,
# per the "updated for chrome and removed recording" commit.

$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("Test Cases")
$wsTestCases.Range("F2").Value = "FAIL"

$wsStartEndDate = $wb.Worksheets.Item("VerifyStartEndDateValidation")
$wsStartEndDate.Range("K2").Value = "FAIL"
